$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.164.18"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.98%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.824.78"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.35%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -1.25%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.37"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.78%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4227"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3672"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.96%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07238"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.66%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.41%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.92"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.37%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.846.64"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.20%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.678"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.91%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07094"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.65%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.292"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.07%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "89.63"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.95%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.27%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008838"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.87%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.003"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.97%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.99"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.278.28"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.56%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.108"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.72%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.88"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.38%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.066.56"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.87%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.975"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.64%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.98"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.54%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.199"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.50%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.39"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.30%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.221"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.16%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.51"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.89%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08815"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.55%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.188"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.30%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.442"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.80%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.002"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.102"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01960"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.76%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05239"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.03%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.238"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.12%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.872"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.02%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1690"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.31%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5034"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.45%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.601"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.46%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.57%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "106.49"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.68%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4735"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.03%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.002"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.13%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06379"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.88%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.36%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.870"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.19%  "

# Row 33/34 swap coin identity (B and C columns) plus updated D/E values
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7451"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.08%  "

$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.956"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.06%  "
